# Daily crypto-price refresh (3-1-2023 -> 4-1-2023) for the "cryptos" sheet:
# updates Price (D), Volume(1h) (E), Data (F) and Hora (G) for every coin row,
# plus a CEJI/KickToken row swap (B/C) at rows 42-43.
# All of D/E/F/G are stored as plain text in the workbook (e.g. "0.05%",
# "4-1-2023", "0"), so NumberFormat is forced to "@" (Text) before each
# assignment - otherwise Excel would auto-convert them to a number, a
# percentage fraction, or a date serial.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.05%"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "4-1-2023"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "0"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.65"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.72%"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "4-1-2023"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "0"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.292"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.45%"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "4-1-2023"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "0"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05707"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.55%"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "4-1-2023"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.636"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.12%"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "4-1-2023"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "0"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.215"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.40%"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "4-1-2023"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8607"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.20%"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "4-1-2023"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8817"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.07%"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "4-1-2023"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1397"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.61%"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "4-1-2023"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07123"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.84%"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "4-1-2023"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03173"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.74%"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "4-1-2023"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09221"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "4-1-2023"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001528"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.23%"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "4-1-2023"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005988"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.06%"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "4-1-2023"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006048"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.70%"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "4-1-2023"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.25%"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "4-1-2023"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "0"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.173"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.95%"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "4-1-2023"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "0"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3125"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.40%"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "4-1-2023"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "0"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.26%"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "4-1-2023"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "0"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.07%"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "4-1-2023"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "0"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.490"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.22%"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "4-1-2023"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04116"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.20%"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "4-1-2023"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "0"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.12%"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = "4-1-2023"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "0"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001220"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.62%"
$ws.Range("F25").NumberFormat = "@"
$ws.Range("F25").Value = "4-1-2023"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "0"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004164"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-16.50%"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "4-1-2023"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "0"

$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "4-1-2023"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "0"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-0.30%"
$ws.Range("F28").NumberFormat = "@"
$ws.Range("F28").Value = "4-1-2023"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "0"

$ws.Range("F29").NumberFormat = "@"
$ws.Range("F29").Value = "4-1-2023"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "0"

$ws.Range("F30").NumberFormat = "@"
$ws.Range("F30").Value = "4-1-2023"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "0"

$ws.Range("F31").NumberFormat = "@"
$ws.Range("F31").Value = "4-1-2023"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "0"

$ws.Range("F32").NumberFormat = "@"
$ws.Range("F32").Value = "4-1-2023"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "0"

$ws.Range("F33").NumberFormat = "@"
$ws.Range("F33").Value = "4-1-2023"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "0"

$ws.Range("F34").NumberFormat = "@"
$ws.Range("F34").Value = "4-1-2023"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "0"

$ws.Range("F35").NumberFormat = "@"
$ws.Range("F35").Value = "4-1-2023"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "0"

$ws.Range("F36").NumberFormat = "@"
$ws.Range("F36").Value = "4-1-2023"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "0"

$ws.Range("F37").NumberFormat = "@"
$ws.Range("F37").Value = "4-1-2023"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "0"

$ws.Range("F38").NumberFormat = "@"
$ws.Range("F38").Value = "4-1-2023"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "0"

$ws.Range("F39").NumberFormat = "@"
$ws.Range("F39").Value = "4-1-2023"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "0"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03801"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.68%"
$ws.Range("F40").NumberFormat = "@"
$ws.Range("F40").Value = "4-1-2023"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "0"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1071"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.05%"
$ws.Range("F41").NumberFormat = "@"
$ws.Range("F41").Value = "4-1-2023"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "0"

$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003770"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-34.47%"
$ws.Range("F42").NumberFormat = "@"
$ws.Range("F42").Value = "4-1-2023"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "0"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002438"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.27%"
$ws.Range("F43").NumberFormat = "@"
$ws.Range("F43").Value = "4-1-2023"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "0"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009483"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.33%"
$ws.Range("F44").NumberFormat = "@"
$ws.Range("F44").Value = "4-1-2023"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "0"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005269"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.49%"
$ws.Range("F45").NumberFormat = "@"
$ws.Range("F45").Value = "4-1-2023"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "0"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.04%"
$ws.Range("F46").NumberFormat = "@"
$ws.Range("F46").Value = "4-1-2023"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "0"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1150"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "101.66%"
$ws.Range("F47").NumberFormat = "@"
$ws.Range("F47").Value = "4-1-2023"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "0"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002262"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.88%"
$ws.Range("F48").NumberFormat = "@"
$ws.Range("F48").Value = "4-1-2023"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "0"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.04%"
$ws.Range("F49").NumberFormat = "@"
$ws.Range("F49").Value = "4-1-2023"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "0"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.04%"
$ws.Range("F50").NumberFormat = "@"
$ws.Range("F50").Value = "4-1-2023"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "0"

$ws.Range("F51").NumberFormat = "@"
$ws.Range("F51").Value = "4-1-2023"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "0"
